# Update cryptocurrency price figures in column D ("Price") on the active
# worksheet. Values are stored as text in the workbook (e.g. "246.80"), so
# each new value is written with a leading apostrophe to force Excel to
# keep it as text instead of silently re-parsing it as a number (which
# would drop meaningful trailing zeros / change the stored representation).
# The style is then reset to "Normal" so no visible number-format change
# (like quote-prefix indicators) lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    "D2"  = "246.80"
    "D3"  = "22.52"
    "D5"  = "0.05685"
    "D6"  = "3.415"
    "D7"  = "6.307"
    "D9"  = "0.8591"
    "D10" = "0.1411"
    "D12" = "0.03023"
    "D13" = "0.03078"
    "D14" = "0.09380"
    "D15" = "3.875"
    "D17" = "0.04763"
    "D18" = "0.0005852"
    "D19" = "0.006406"
    "D20" = "0.005029"
    "D21" = "0.0009967"
    "D22" = "0.0001500"
    "D24" = "2.194"
    "D26" = "0.1283"
    "D41" = "0.006830"
    "D42" = "0.1065"
    "D43" = "0.002670"
    "D44" = "0.008436"
    "D46" = "0.00000000750"
    "D47" = "0.4501"
    "D48" = "0.1476"
    "D50" = "0.01010"
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    $cell.Value = "'" + $updates[$addr]
    $cell.Style = "Normal"
}
